# Add "status_label" as string version of "status".
#
# This inserts a new column B ("status_label") before the existing "NCTId"
# column, shifting every column from B..I to C..J. It then fills the new
# column with the French label corresponding to the existing emoji
# "statut" (status) column A ("🟥" => "rouge") for every data row.
#
# Additionally, the underlying data source used to generate this
# publipostage sheet was re-ordered for rows 3 and 4 (the NCT03433859 /
# SALUTOX trial and the NCT04870814 tourniquet trial swapped places), so
# those two rows' NCTId / clinical_trial_title / acronym values are
# swapped back into the row order seen in the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new "status_label" column at B, shifting B:I -> C:J ---
$ws.Columns.Item(2).Insert()

$ws.Cells.Item(1,2).Value = "status_label"

# Map the existing emoji "statut" values (column A) to their French text
# label and write them into the new column B for every data row.
$statusLabels = @{
    "🟥" = "rouge"
    "🟩" = "vert"
    "🟨" = "jaune"
    "🟧" = "orange"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $status = $ws.Cells.Item($r, 1).Value()
    $label = $statusLabels[$status]
    if (-not $label) { $label = $status }
    $ws.Cells.Item($r, 2).Value = $label
}

# --- 2. Swap the NCTId / clinical_trial_title / acronym data that moved
#        between rows 3 and 4 (columns C, F and G after the insert) ---
$c3 = $ws.Cells.Item(3,3).Value()
$f3 = $ws.Cells.Item(3,6).Value()
$g3 = $ws.Cells.Item(3,7).Value()

$c4 = $ws.Cells.Item(4,3).Value()
$f4 = $ws.Cells.Item(4,6).Value()
$g4 = $ws.Cells.Item(4,7).Value()

$ws.Cells.Item(3,3).Value = $c4
$ws.Cells.Item(3,6).Value = $f4
if ($g4) { $ws.Cells.Item(3,7).Value = $g4 } else { $ws.Cells.Item(3,7).Value = "" }

$ws.Cells.Item(4,3).Value = $c3
$ws.Cells.Item(4,6).Value = $f3
if ($g3) { $ws.Cells.Item(4,7).Value = $g3 } else { $ws.Cells.Item(4,7).Value = "" }
